{"js": "// The document contains three \"<id>p138v_N</id>\" markers, each currently\n// split across three separate runs:\n//   run1: \"<id>\"      (Courier New, color 7f6000, sz 18)\n//   run2: \"p138v_N\"   (plain, color 000000)\n//   run3: \"</id>\"     (Courier New, color 7f6000, sz 18)\n// The edit merges each triplet into a single run whose text is the\n// concatenation \"<id>p138v_N</id>\" and whose formatting matches the\n// surrounding Courier New / 7f6000 / sz18 runs (i.e. the formatting of the\n// first run in the triplet).\nconst body = context.document.body;\n\nfor (const n of [1, 2, 3]) {\n  const needle = `<id>p138v_${n}</id>`;\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    // Already merged (or not present) - nothing to do for this marker.\n    continue;\n  }\n\n  // Replacing the whole matched range with the same text collapses the\n  // three backing runs into a single run that inherits the formatting of\n  // the range's first run (\"<id>\", Courier New / 7f6000 / sz 18).\n  const found = results.items[0];\n  found.insertText(needle, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The document contains three \"<id>p138v_N</id>\" markers, each currently\n# split across three separate runs:\n#   run1: \"<id>\"      (Courier New, color 7f6000, sz 18)\n#   run2: \"p138v_N\"   (plain, color 000000)\n#   run3: \"</id>\"     (Courier New, color 7f6000, sz 18)\n# The edit merges each triplet into a single run whose text is the\n# concatenation \"<id>p138v_N</id>\" and whose formatting matches the\n# surrounding Courier New / 7f6000 / sz18 runs (i.e. the formatting of the\n# first run in the triplet). A Find/Replace over the whole document content\n# collapses the three backing runs into one run that inherits the\n# formatting of the first (\"<id>\") run.\n$d = $word.ActiveDocument\n\nforeach ($n in 1..3) {\n    $needle = \"<id>p138v_$n</id>\"\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2)\n    Write-Output \"p138v_$n replaced: $found\"\n}\n"}
